# Weekly update: insert a new price record for "Ajo / Chino / Primera" at
# Terminal Hortofrutícola Agro Chillán. The new record is inserted right
# before the existing row 148, pushing rows 148-207 down to 149-208 (the
# sheet's A1:R207 range becomes A1:R208). All other existing rows/columns
# (A,B,C + everything outside this block) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 148 - shifts rows 148:207 down to 149:208 and
# carries the column D (date) number format down onto the new row.
$ws.Rows.Item(148).Insert()

# Populate the new row 148 with the new weekly record.
$ws.Cells.Item(148, 1).Value = 7
$ws.Cells.Item(148, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(148, 3).Value = "Ñuble"
$ws.Cells.Item(148, 4).Value = 44636
$ws.Cells.Item(148, 5).Value = 16
$ws.Cells.Item(148, 6).Value = 100112003
$ws.Cells.Item(148, 7).Value = "Ajo"
$ws.Cells.Item(148, 8).Value = "Chino"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 60
$ws.Cells.Item(148, 11).Value = 19000
$ws.Cells.Item(148, 12).Value = 20000
$ws.Cells.Item(148, 13).Value = 19500
$ws.Cells.Item(148, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(148, 15).Value = "China"
$ws.Cells.Item(148, 16).Value = 1950
$ws.Cells.Item(148, 17).Value = 10
$ws.Cells.Item(148, 18).Value = "Hortaliza"
